$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 171.831087874
$ws.Range("C2").Value = 21.0028877056

$ws.Range("B3").Value = 171831.087874
$ws.Range("C3").Value = 21002.8877056

$ws.Range("B4").Value = 343662.175748
$ws.Range("C4").Value = 63008.66311679999

$ws.Range("B5").Value = 6873.24351496
$ws.Range("C5").Value = 1260.173262336
